# Daily TGP (terminal gate pricing) refresh.
# Rolls the "Effective Date" pricing table forward one day:
#   - rows previously dated 6 Feb 2026 (serial 46059) now show the prior
#     "7 Feb 2026" (46060) figures,
#   - rows previously dated 7 Feb 2026 (serial 46060) now show brand-new
#     10 Feb 2026 (46063) prices.
# Terminal names/labels and cell styles are untouched; only the date and
# Diesel/ULP/PULP/e10 price columns move.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A8").Value = 46063
$ws.Range("D8").Value = 159.35
$ws.Range("E8").Value = 148.62
$ws.Range("F8").Value = 158.62
$ws.Range("G8").Value = 148.51

$ws.Range("A9").Value = 46063
$ws.Range("D9").Value = 159.35
$ws.Range("E9").Value = 148.62
$ws.Range("F9").Value = 158.62
$ws.Range("G9").Value = 148.51

$ws.Range("A10").Value = 46063
$ws.Range("D10").Value = 160.74
$ws.Range("E10").Value = 150.77
$ws.Range("F10").Value = 160.77
$ws.Range("G10").Value = 151.02

$ws.Range("A11").Value = 46060
$ws.Range("D11").Value = 159.14
$ws.Range("E11").Value = 148.6
$ws.Range("F11").Value = 158.6
$ws.Range("G11").Value = 148.49

$ws.Range("A12").Value = 46060
$ws.Range("D12").Value = 159.14
$ws.Range("E12").Value = 148.6
$ws.Range("F12").Value = 158.6
$ws.Range("G12").Value = 148.49

$ws.Range("A13").Value = 46060
$ws.Range("D13").Value = 160.48
$ws.Range("E13").Value = 150.71
$ws.Range("F13").Value = 160.71
$ws.Range("G13").Value = 150.95

$ws.Range("A17").Value = 46063
$ws.Range("D17").Value = 164.95
$ws.Range("E17").Value = 154.17
$ws.Range("F17").Value = 164.17

$ws.Range("A18").Value = 46060
$ws.Range("D18").Value = 164.68
$ws.Range("E18").Value = 154.09
$ws.Range("F18").Value = 164.09

$ws.Range("A22").Value = 46063
$ws.Range("D22").Value = 160.42
$ws.Range("E22").Value = 150.74
$ws.Range("F22").Value = 160.34
$ws.Range("G22").Value = 152.49

$ws.Range("A23").Value = 46063
$ws.Range("D23").Value = 165.73
$ws.Range("E23").Value = 156.45
$ws.Range("F23").Value = 166.45
$ws.Range("G23").Value = "N/A"

$ws.Range("A24").Value = 46063
$ws.Range("D24").Value = 165.89
$ws.Range("E24").Value = 157.06
$ws.Range("F24").Value = 167.06
$ws.Range("G24").Value = "N/A"

$ws.Range("A25").Value = 46063
$ws.Range("D25").Value = 165.89
$ws.Range("E25").Value = 156.59
$ws.Range("F25").Value = 166.59
$ws.Range("G25").Value = 157.44

$ws.Range("A26").Value = 46063
$ws.Range("D26").Value = 165.5
$ws.Range("E26").Value = 158.17
$ws.Range("F26").Value = 168.17
$ws.Range("G26").Value = "N/A"

$ws.Range("A27").Value = 46060
$ws.Range("D27").Value = 160.21
$ws.Range("E27").Value = 150.72
$ws.Range("F27").Value = 160.32
$ws.Range("G27").Value = 152.47

$ws.Range("A28").Value = 46060
$ws.Range("D28").Value = 165.47
$ws.Range("E28").Value = 156.39
$ws.Range("F28").Value = 166.39
$ws.Range("G28").Value = "N/A"

$ws.Range("A29").Value = 46060
$ws.Range("D29").Value = 165.64
$ws.Range("E29").Value = 157.01
$ws.Range("F29").Value = 167.01
$ws.Range("G29").Value = "N/A"

$ws.Range("A30").Value = 46060
$ws.Range("D30").Value = 165.63
$ws.Range("E30").Value = 156.53
$ws.Range("F30").Value = 166.53
$ws.Range("G30").Value = 157.39

$ws.Range("A31").Value = 46060
$ws.Range("D31").Value = 165.24
$ws.Range("E31").Value = 158.12
$ws.Range("F31").Value = 168.12
$ws.Range("G31").Value = "N/A"

$ws.Range("A35").Value = 46063
$ws.Range("D35").Value = 159.23
$ws.Range("E35").Value = 148.46
$ws.Range("F35").Value = 157.46

$ws.Range("A36").Value = 46060
$ws.Range("D36").Value = 158.97
$ws.Range("E36").Value = 148.4
$ws.Range("F36").Value = 157.4

$ws.Range("A40").Value = 46063
$ws.Range("D40").Value = 165.56
$ws.Range("E40").Value = 156.19
$ws.Range("F40").Value = 166.19

$ws.Range("A41").Value = 46063
$ws.Range("D41").Value = 165.27
$ws.Range("E41").Value = 156.61
$ws.Range("F41").Value = 166.61

$ws.Range("A42").Value = 46060
$ws.Range("D42").Value = 165.31
$ws.Range("E42").Value = 156.08
$ws.Range("F42").Value = 166.08

$ws.Range("A43").Value = 46060
$ws.Range("D43").Value = 165.03
$ws.Range("E43").Value = 156.51
$ws.Range("F43").Value = 166.51

$ws.Range("A47").Value = 46063
$ws.Range("D47").Value = 160.19
$ws.Range("E47").Value = 150.17
$ws.Range("F47").Value = 160.17

$ws.Range("A48").Value = 46063
$ws.Range("D48").Value = 159.85
$ws.Range("E48").Value = 150.13
$ws.Range("F48").Value = 160.13

$ws.Range("A49").Value = 46060
$ws.Range("D49").Value = 159.74
$ws.Range("E49").Value = 150.03
$ws.Range("F49").Value = 160.03

$ws.Range("A50").Value = 46060
$ws.Range("D50").Value = 159.39
$ws.Range("E50").Value = 149.98
$ws.Range("F50").Value = 159.98

$ws.Range("A54").Value = 46063
$ws.Range("D54").Value = 174.64
$ws.Range("E54").Value = 163.95
$ws.Range("F54").Value = 173.95

$ws.Range("A55").Value = 46063
$ws.Range("D55").Value = 164.02
$ws.Range("E55").Value = 161.88
$ws.Range("F55").Value = 171.88

$ws.Range("A56").Value = 46063
$ws.Range("D56").Value = 163.89
$ws.Range("E56").Value = "N/A"
$ws.Range("F56").Value = "N/A"

$ws.Range("A57").Value = 46063
$ws.Range("D57").Value = 164.56
$ws.Range("E57").Value = 156.3
$ws.Range("F57").Value = "N/A"

$ws.Range("A58").Value = 46063
$ws.Range("D58").Value = 160.33
$ws.Range("E58").Value = 152.2
$ws.Range("F58").Value = 162.2

$ws.Range("A59").Value = 46063
$ws.Range("D59").Value = 167.34
$ws.Range("E59").Value = 162.25
$ws.Range("F59").Value = "N/A"

$ws.Range("A60").Value = 46060
$ws.Range("D60").Value = 174.37
$ws.Range("E60").Value = 163.95
$ws.Range("F60").Value = 173.95

$ws.Range("A61").Value = 46060
$ws.Range("D61").Value = 163.76
$ws.Range("E61").Value = 161.8
$ws.Range("F61").Value = 171.8

$ws.Range("A62").Value = 46060
$ws.Range("D62").Value = 163.63
$ws.Range("E62").Value = "N/A"
$ws.Range("F62").Value = "N/A"

$ws.Range("A63").Value = 46060
$ws.Range("D63").Value = 164.28
$ws.Range("E63").Value = 156.22
$ws.Range("F63").Value = "N/A"

$ws.Range("A64").Value = 46060
$ws.Range("D64").Value = 160.05
$ws.Range("E64").Value = 152.12
$ws.Range("F64").Value = 162.12

$ws.Range("A65").Value = 46060
$ws.Range("D65").Value = 167.07
$ws.Range("E65").Value = 162.23
$ws.Range("F65").Value = "N/A"
